$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Row 3: GTX gain settings used for the three measurement columns ---
$ws.Range("B3").Value = 3.0
$ws.Range("C3").Value = 15.0
$ws.Range("D3").Value = 30.0

# --- Row 4: series labels, left block (B:D) and the mirrored right block (G:I) ---
$ws.Range("B4").Value = "GTX = 3 dB"
$ws.Range("C4").Value = "GTX = 15 dB"
$ws.Range("D4").Value = "GTX = 30 dB"
$ws.Range("G4").Value = "GTX = 3 dB"
$ws.Range("H4").Value = "GTX = 15 dB"
$ws.Range("I4").Value = "GTX = 30 dB"

# --- Measured power data (columns B, C, D) for rows 5-19, in frequency order ---
# Each row: row number, then the B/C/D measured values.
$rows = @(
    , @(5,  -44.75, -32.98, -18.59)
    , @(6,  -43.54, -31.83, -17.42)
    , @(7,  -43.22, -31.46, -17.00)
    , @(8,  -43.19, -31.41, -16.90)
    , @(9,  -43.35, -31.60, -17.05)
    , @(10, -43.52, -31.76, -17.20)
    , @(11, -45.60, -33.88, -19.33)
    , @(12, -48.02, -36.44, -21.92)
    , @(13, -50.90, -39.43, -24.93)
    , @(14, -53.03, -41.63, -27.15)
    , @(15, -53.42, -41.91, -27.44)
    , @(16, -54.64, -43.10, -28.64)
    , @(17, -55.81, -44.20, -29.67)
    , @(18, -57.83, -46.23, -31.69)
    , @(19, -59.31, -47.66, -33.08)
)

foreach ($entry in $rows) {
    $r = $entry[0]
    $ws.Cells.Item($r, 2).Value = $entry[1]
    $ws.Cells.Item($r, 3).Value = $entry[2]
    $ws.Cells.Item($r, 4).Value = $entry[3]
}

# Rows 6-19 (B:D) pick up the same 2-decimal number-format style row 5 already has.
$ws.Range("B5:D5").Copy()
$ws.Range("B6:D19").PasteSpecial(-4122)

# --- Reference transmitter power level used by the G/H/I attenuation formulas ---
$ws.Range("B21").Value = -16.16
